$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.482.45'
$ws.Range('E2').Value = '  -1.75%  '
$ws.Range('D3').Value = '2.336.39'
$ws.Range('E3').Value = '  -4.14%  '
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '541.75'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.58'
$ws.Range('E6').Value = '  -6.63%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.520'
$ws.Range('E8').Value = '  -10.29%  '
$ws.Range('D9').Value = '2.334.96'
$ws.Range('E9').Value = '  -4.45%  '
$ws.Range('E10').Value = '  -1.77%  '
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.27'
$ws.Range('E12').Value = '  -3.18%  '
$ws.Range('E13').Value = '  -2.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '24.38'
$ws.Range('E14').Value = '  -5.33%  '
$ws.Range('D15').Value = '2.749.94'
$ws.Range('D16').Value = '60.170.29'
$ws.Range('E16').Value = '  -2.62%  '
$ws.Range('E17').Value = '  -4.19%  '
$ws.Range('D18').Value = '2.331.63'
$ws.Range('E18').Value = '  -4.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.53'
$ws.Range('E19').Value = '  -1.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '317.10'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.06'
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.51'
$ws.Range('E22').Value = '  -5.17%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.98'
$ws.Range('E24').Value = '  -1.41%  '
$ws.Range('E25').Value = '  -8.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.29'
$ws.Range('E26').Value = '  +6.39%  '
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.85'
$ws.Range('E28').Value = '  -3.44%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '493.04'
$ws.Range('E29').Value = '  -5.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.36'
$ws.Range('E30').Value = '  -7.24%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0850'
$ws.Range('E31').Value = '  -11.63%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.144'
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.77'
$ws.Range('E33').Value = '  -4.79%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.49'
$ws.Range('E34').Value = '  -5.66%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.995'
$ws.Range('E35').Value = '  -0.46%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.53'
$ws.Range('E36').Value = '  -3.75%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.44'
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.371'
$ws.Range('E38').Value = '  -1.95%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.14'
$ws.Range('E39').Value = '  -8.18%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.76'
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '140.87'
$ws.Range('E41').Value = '  +1.70%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.51'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '140.29'
$ws.Range('E44').Value = '  -1.21%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.52'
$ws.Range('E45').Value = '  -1.31%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.03'
$ws.Range('E46').Value = '  -9.81%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0507'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.89'
$ws.Range('E48').Value = '  -10.20%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.564'
$ws.Range('E49').Value = '  -3.99%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0894'
$ws.Range('E50').Value = '  -3.35%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0218'
$ws.Range('E51').Value = '  -3.72%  '
